$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "B"
$ws.Range("E3").Value = "B"

$ws.Range("E7").Select()
